$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "[Inches]" -> "[" + "Foaming " + "Inches]" (three runs,
# same run formatting as the original run).
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("[Inches]", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Shrink the found range down to just "[" and leave the cursor
    # collapsed right after it.
    $rng.Text = "["
    $rng.Collapse(0)

    # Insert "Foaming " right after "[". Briefly toggle Bold on/off so
    # this text is not silently re-coalesced into the previous run even
    # though the final formatting is identical.
    $r2 = $rng.Duplicate
    $r2.InsertAfter("Foaming ")
    $r2.Bold = 1

    # Insert "Inches]" right after "Foaming ". Same trick: toggle Bold
    # on the newly inserted text to keep it a distinct run, using a
    # fresh collapsed range so $r2's span is not disturbed.
    $ins = $r2.Duplicate
    $ins.Collapse(0)
    $r3 = $ins.Duplicate
    $r3.InsertAfter("Inches]")
    $r3.Bold = 1
    $r3.Bold = 0

    # Now restore $r2 ("Foaming ") back to non-bold too.
    $r2.Bold = 0
}

# ---------------------------------------------------------------------
# Change 2: merge the "  " run + "[TotalInches]" run into a single run
# "  [TotalInches]" (keeps the first run's formatting, which already
# matches the second run's formatting in the source document).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("  [TotalInches]", $false, $false, $false, $false, $false, $true, 1, $false, "  [TotalInches]", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3: merge the "  " run + "[Dimensions]" run into a single run
# "  [Dimensions]".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("  [Dimensions]", $false, $false, $false, $false, $false, $true, 1, $false, "  [Dimensions]", 2) | Out-Null
